# update base, test signUp browserStack
# Refresh the phone-number / OTP test data used by the SignUp flow
# (Login sheet's OTP phone number + SignUp sheet's phone & OTP values).

$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("Login")
$wsSignUp = $wb.Worksheets.Item("SignUp")

# Login sheet: OTP test-case phone number.
# Leading apostrophe forces text entry so the leading zero survives even
# though the cell's existing number format is numeric.
$wsLogin.Range("B12").Value = "'0363714939"

# SignUp sheet: phone number used for registration + the OTP values sent
# for it (duplicated across the "enter"/"resend" rows).
$wsSignUp.Range("B2").Value = "0363870101"
$wsSignUp.Range("E4").Value = "922734"
$wsSignUp.Range("E5").Value = "922734"
$wsSignUp.Range("E6").Value = "752110"
$wsSignUp.Range("E7").Value = "752110"
$wsSignUp.Range("E8").Value = "958279"
